# Generate Report for Handoff
# Updates the localization-status workbook: the 52080a9f... file finished
# handoff, so its Priority flips from "low" -> "ht" and its Latest Handoff
# Datetime advances, for both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 (all rows tied to the 52080a9f... source file)
$ws_zhcn.Range("E4:E7").Value = "ht"
$ws_zhcn.Range("H4:H7").Value = "2016-08-21 10:39:28"

# de-de: rows 4-7 (same source file, its own language-specific handoff xlf)
$ws_dede.Range("E4:E7").Value = "ht"
$ws_dede.Range("H4:H7").Value = "2016-08-21 10:39:32"

# Overview: "Latest HO Xliff Generate Date" mirrors the de-de handoff time
$ws_overview = $wb.Worksheets.Item("Overview")
$ws_overview.Range("G4:G7").Value = "2016-08-21 10:39:32"
